$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the A1 header text (1.5.4 -> 13.1.3 indicator code, rest unchanged)
$ws.Range("A1").Value = "13.1.3 Кырсыктардын кооптуулугун азайтуунун улуттук стратегияларына ылайык, кырсыктардын кооптуулугун азайтуунун жергиликтүү стратегияларын кабыл алган жана ишке ашырган жергиликтүү бийлик органдарынын үлүшү"

# 2. Extend the table with columns E:H (years 2020-2023), copying the
#    number formatting/style from column D so the new cells match.
$ws.Range("D3:D6").Copy($ws.Range("E3"))
$ws.Range("D3:D6").Copy($ws.Range("F3"))
$ws.Range("D3:D6").Copy($ws.Range("G3"))
$ws.Range("D3:D6").Copy($ws.Range("H3"))

# Row 3: year headers
$ws.Range("D3").Value = 2019
$ws.Range("E3").Value = 2020
$ws.Range("F3").Value = 2021
$ws.Range("G3").Value = 2022
$ws.Range("H3").Value = 2023

# Row 4: number of local governments (484, now numeric instead of text)
$ws.Range("D4").Value = 484
$ws.Range("E4").Value = 484
$ws.Range("F4").Value = 484
$ws.Range("G4").Value = 484
$ws.Range("H4").Value = 484

# Row 5: proportion (%)
$ws.Range("E5").Value = 13.2
$ws.Range("F5").Value = 21.5
$ws.Range("G5").Value = 34.5
$ws.Range("H5").Value = 40.53

# Row 6: count of governments implementing DRR strategies
$ws.Range("E6").Value = 67
$ws.Range("F6").Value = 104
$ws.Range("G6").Value = 167
$ws.Range("H6").Value = 169
